# Mod Part Patching for 0.8.0
#
# The SSS Category / Tier columns were reworked (categories renamed /
# split, tiers added) and the table was then re-sorted by Category then
# Tier via Data > Sort. This rebuilds the table in its final, sorted
# form (including the per-row TechRequired formula, now unshared and
# without the old trailing blank line) and finally re-applies an
# equivalent sort so the workbook records the same <sortState>.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime the shared-string table so the three brand-new category labels
# get interned in the same order the source workbook uses them in:
# landing, comms, reactors (the loop below re-asserts the real values).
$ws.Cells.Item(26, 2).Value2 = "landing"
$ws.Cells.Item(2, 2).Value2 = "comms"
$ws.Cells.Item(29, 2).Value2 = "reactors"

# target row | PART_name | SSS Category | Tier (blank = no tier)
$rowData = @"
2|bluedog_RAE_TelemetryAntenna|comms|4
3|bluedog_Pioneer_HGA_Alternate|comms|6
4|bluedog_Pioneer_HGA|comms|6
5|bluedog_Pioneer_LGA|comms|6
6|bluedog_Pioneer_ProbeAntenna|comms|6
7|bluedog_Pioneer_RelayAntenna|comms|6
8|bluedog_RAE_MotorDecoupler|construction|4
9|bluedog_Titan2_Adapter|construction|5
10|bluedog_Pioneer_AftAdapter|construction|6
11|bluedog_Pioneer_AftAdapterAlternate|construction|6
12|bluedog_Pioneer_Decoupler|construction|6
13|bluedog_Pioneer_ExperimentBus|construction|6
14|bluedog_Pioneer_ExtensionBus|construction|6
15|bluedog_Pioneer_ProbeAdapter|construction|6
16|bluedog_Pioneer_ProbeDecoupler|construction|6
17|bluedog_RAE_RCS|control|4
18|bluedog_Pioneer_StellarReferenceAssembly|control|6
19|bluedog_LOantenna|debug|
20|bluedog_LOdish|debug|
21|bluedog_upgrade_logo|debug|
22|bluedog_solarBattery|electrics|
23|bluedog_RAE_VCPS|hypergol|4
24|bluedog_Pioneer_OrbiterPropulsionUnit|hypergol|6
25|bluedog_SmallKlaw|isru|9
26|bluedog_Pioneer_Probe|landing|6
27|bluedog_RAE_ProbeCore|probes|4
28|bluedog_Pioneer_ProbeCore|probes|6
29|bluedog_Pioneer_SNAP19|reactors|6
30|bluedog_Explorer_Beacon_SLR|science|3
31|bluedog_RAE_DipoleAntenna|science|4
32|bluedog_RAE_RadioAntenna|science|4
33|bluedog_Pioneer_AMD|science|6
34|bluedog_Pioneer_CPD|science|6
35|bluedog_Pioneer_CRT|science|6
36|bluedog_Pioneer_GeigerTube|science|6
37|bluedog_Pioneer_InfraredRadiometer|science|6
38|bluedog_Pioneer_IPP|science|6
39|bluedog_Pioneer_Magnetometer|science|6
40|bluedog_Pioneer_MLS|science|6
41|bluedog_Pioneer_OrbiterScanner|science|6
42|bluedog_Pioneer_PlasmaAnalyzer|science|6
43|bluedog_Pioneer_ProbeExperiment|science|6
44|bluedog_Pioneer_TRD|science|6
45|bluedog_Pioneer_UltravioletPhotometer|science|6
46|bluedog_RAE_SolarPaddle|solar|4
47|bluedog_RAE_TrackingSolarPaddle|solar|4
48|bluedog_Explorer_Star17|solids|3
49|bluedog_Pioneer_LargeLouver|thermal|6
50|bluedog_Pioneer_MediumLouver|thermal|6
51|bluedog_Pioneer_Radiator|thermal|6
52|bluedog_Pioneer_SmallLouver|thermal|6
"@

$lines = $rowData -split "`n" | Where-Object { $_.Trim() -ne "" }

foreach ($line in $lines) {
    $parts = $line.Trim() -split '\|'
    $r = [int]$parts[0]
    $partName = $parts[1]
    $cat = $parts[2]
    $tierText = $parts[3]

    # Column A: PART_name
    $ws.Cells.Item($r, 1).Value2 = $partName

    # Column B: SSS Category
    $ws.Cells.Item($r, 2).Value2 = $cat

    # Column C: Tier (left blank when not applicable)
    if ([string]::IsNullOrEmpty($tierText)) {
        $ws.Cells.Item($r, 3).Value2 = $null
    } else {
        $ws.Cells.Item($r, 3).Value2 = [int]$tierText
    }

    # Column D: Mod (unchanged, but set for safety/consistency)
    $ws.Cells.Item($r, 4).Value2 = "Bluedog_DB"

    # Column E: rebuild the (now unshared) formula for this row without
    # the extra trailing blank line that the old shared formula had.
    $formula = "=""@PART[""&A$r&""]:AFTER[""&D$r&""] //`n{`n`t@TechRequired = ""&B$r&C$r&""`n}"""
    $ws.Cells.Item($r, 5).Formula = $formula

    # Setting a multi-line formula can make the COM layer auto-grow the
    # row height; AutoFit it back down so no stray row height survives.
    $ws.Rows.Item($r).AutoFit() | Out-Null
}

# Re-apply Data > Sort over the table (A2:E52) by Category (B) then
# Tier (C), ascending. The data above is already arranged in this
# order, so this reproduces the workbook's <sortState> without moving
# anything.
$sortRange = $ws.Range("A2:E52")
$key1 = $ws.Range("B2:B52")
$key2 = $ws.Range("C2:C52")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($key1)
$ws.Sort.SortFields.Add($key2)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

$wb.Save()
